# Applies the "Add files via upload" update to maine_covid_summary.xlsx
#  - cases_by_race: renumber the 2020-12-16 block's index column and
#    append a new 2020-12-20 / 2020-12-19 refresh block (rows 71-79)
#  - cases_by_ethnicity: renumber the 2020-12-16 block's index column and
#    append a new 2020-12-20 / 2020-12-19 refresh block (rows 29-31)

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    if ($text -eq "") {
        $cell.Value = ""
    } else {
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------
# Sheet 1: cases_by_race
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("cases_by_race")

# Renumber the running index for the existing 2020-12-16 block
# (rows 62-70 held 0..8, they continue the running count as 60..68)
$raceReindex = @{ 62 = 60; 63 = 61; 64 = 62; 65 = 63; 66 = 64; 67 = 65; 68 = 66; 69 = 67; 70 = 68 }
foreach ($r in $raceReindex.Keys) {
    $ws1.Cells.Item($r, 1).Value = $raceReindex[$r]
}

# Append the new 2020-12-20 (as-of 2020-12-19) block
$raceRows = @(
    @{ Row = 71; Idx = 0; Label = ""; Cases = 7 },
    @{ Row = 72; Idx = 1; Label = "American Indian or Alaska Native"; Cases = 59 },
    @{ Row = 73; Idx = 2; Label = "Asian"; Cases = 252 },
    @{ Row = 74; Idx = 3; Label = "Black or African American"; Cases = 1465 },
    @{ Row = 75; Idx = 4; Label = "Native Hawaiian or Other Pacific Islander"; Cases = 10 },
    @{ Row = 76; Idx = 5; Label = "Not disclosed"; Cases = 1925 },
    @{ Row = 77; Idx = 6; Label = "Other Race"; Cases = 418 },
    @{ Row = 78; Idx = 7; Label = "Two or more"; Cases = 158 },
    @{ Row = 79; Idx = 8; Label = "White"; Cases = 14652 }
)

foreach ($entry in $raceRows) {
    $r = $entry.Row
    # Column A carries the bold/bordered "index" style (s="1") throughout
    # the sheet; copy it from the row above (formats only get copied,
    # then the real value is written over it).
    $ws1.Cells.Item($r - 1, 1).Copy($ws1.Cells.Item($r, 1))
    $ws1.Cells.Item($r, 1).Value = $entry.Idx
    Set-TextCell $ws1 $r 2 $entry.Label
    Set-TextCell $ws1 $r 3 "2020-12-20"
    Set-TextCell $ws1 $r 4 "2020-12-19"
    $ws1.Cells.Item($r, 5).Value = $entry.Cases
}

# ---------------------------------------------------------------
# Sheet 2: cases_by_ethnicity
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("cases_by_ethnicity")

# Renumber the running index for the existing 2020-12-16 block
# (rows 26-28 held 0..2, they continue the running count as 24..26)
$ethReindex = @{ 26 = 24; 27 = 25; 28 = 26 }
foreach ($r in $ethReindex.Keys) {
    $ws2.Cells.Item($r, 1).Value = $ethReindex[$r]
}

# Append the new 2020-12-20 (as-of 2020-12-19) block
$ethRows = @(
    @{ Row = 29; Idx = 0; Label = "unknown"; Cases = 4140 },
    @{ Row = 30; Idx = 1; Label = "Not Hispanic or Latino"; Cases = 14416 },
    @{ Row = 31; Idx = 2; Label = "Hispanic or Latino"; Cases = 390 }
)

foreach ($entry in $ethRows) {
    $r = $entry.Row
    $ws2.Cells.Item($r - 1, 1).Copy($ws2.Cells.Item($r, 1))
    $ws2.Cells.Item($r, 1).Value = $entry.Idx
    Set-TextCell $ws2 $r 2 $entry.Label
    Set-TextCell $ws2 $r 3 "2020-12-20"
    Set-TextCell $ws2 $r 4 "2020-12-19"
    $ws2.Cells.Item($r, 5).Value = $entry.Cases
}
